# Update receptor (Fzd2) expression-derived metrics with new TPM values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 0.513343
$ws.Range("N2").Value = 1.540029
$ws.Range("O2").Value = 0.02896400434489499
$ws.Range("P2").Value = 0.02896400434489499
$ws.Range("Q2").Value = 0.08208371681433334
$ws.Range("R2").Value = 0.7387534513289999
$ws.Range("S2").Value = 0.000766900780748904
$ws.Range("T2").Value = 0.000766900780748904

# Row 3
$ws.Range("O3").Value = 0.8830650561604291
$ws.Range("P3").Value = 0.8830650561604291
$ws.Range("S3").Value = 0.0233815487995834
$ws.Range("T3").Value = 0.0233815487995834

# Row 4
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.04863733333333333
$ws.Range("N4").Value = 0.145912
$ws.Range("O4").Value = 0.00274423131121058
$ws.Range("P4").Value = 0.00274423131121058
$ws.Range("Q4").Value = 0.007777125812444444
$ws.Range("R4").Value = 0.069994132312
$ws.Range("S4").Value = 0.00007266098672209035
$ws.Range("T4").Value = 0.00007266098672209036

# Row 5
$ws.Range("M5").Value = 1.499596
$ws.Range("N5").Value = 4.498788
$ws.Range("O5").Value = 0.08461068926543686
$ws.Range("P5").Value = 0.08461068926543686
$ws.Range("Q5").Value = 0.2397859002653334
$ws.Range("R5").Value = 2.158073102388
$ws.Range("S5").Value = 0.0022402980915449
$ws.Range("T5").Value = 0.0022402980915449

# Row 6
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.010918
$ws.Range("N6").Value = 0.032754
$ws.Range("O6").Value = 0.0006160189180286153
$ws.Range("P6").Value = 0.0006160189180286153
$ws.Range("Q6").Value = 0.001745791839333333
$ws.Range("R6").Value = 0.015712126554
$ws.Range("S6").Value = 0.00001631077607801515
$ws.Range("T6").Value = 0.00001631077607801516

# Row 7
$ws.Range("M7").Value = 0.513343
$ws.Range("N7").Value = 1.540029
$ws.Range("O7").Value = 0.02896400434489499
$ws.Range("P7").Value = 0.02896400434489499
$ws.Range("Q7").Value = 3.018021525136
$ws.Range("R7").Value = 27.162193726224
$ws.Range("S7").Value = 0.02819710356414609
$ws.Range("T7").Value = 0.02819710356414609

# Row 8
$ws.Range("O8").Value = 0.8830650561604291
$ws.Range("P8").Value = 0.8830650561604291
$ws.Range("S8").Value = 0.8596835073608456
$ws.Range("T8").Value = 0.8596835073608458

# Row 9
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.04863733333333333
$ws.Range("N9").Value = 0.145912
$ws.Range("O9").Value = 0.00274423131121058
$ws.Range("P9").Value = 0.00274423131121058
$ws.Range("Q9").Value = 0.2859462755413333
$ws.Range("R9").Value = 2.573516479872
$ws.Range("S9").Value = 0.002671570324488489
$ws.Range("T9").Value = 0.002671570324488489

# Row 10
$ws.Range("M10").Value = 1.499596
$ws.Range("N10").Value = 4.498788
$ws.Range("O10").Value = 0.08461068926543686
$ws.Range("P10").Value = 0.08461068926543686
$ws.Range("Q10").Value = 8.816352822592002
$ws.Range("R10").Value = 79.34717540332801
$ws.Range("S10").Value = 0.08237039117389196
$ws.Range("T10").Value = 0.08237039117389196

# Row 11
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.010918
$ws.Range("N11").Value = 0.032754
$ws.Range("O11").Value = 0.0006160189180286153
$ws.Range("P11").Value = 0.0006160189180286153
$ws.Range("Q11").Value = 0.064188581536
$ws.Range("R11").Value = 0.577697233824
$ws.Range("S11").Value = 0.0005997081419506002
$ws.Range("T11").Value = 0.0005997081419506003

Write-Host "Applied updated TPM-derived values to 80 cells."
